$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the pre-existing "_GoBack" bookmark that currently sits at
#    the very start of the document (around "D I O C E S E ...").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Replace the signer's name with a blank line for the final
#    recommendation / signature ("Rev. Fr. Ruben C. Espinosa" ->
#    "_______________________").
# ------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute("Rev. Fr. Ruben C. Espinosa", $true, $false, $false, $false, $false, $true, 1, $false, `
    "_______________________", 2)

# ------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark right after the new run, i.e. at
#    the end of that paragraph (immediately before its paragraph
#    mark) -- mirroring where Word leaves it after the last edit.
#
#    A zero-width Range sitting exactly on "end of last run in a
#    paragraph" cannot be targeted directly, so a temporary marker
#    character is inserted right after the text, the bookmark is
#    placed in the (now perfectly safe, mid-paragraph) gap in front
#    of it, and the marker is removed again -- leaving the bookmark
#    collapsed right after the underscores, before the paragraph
#    mark, exactly like the target edit.
# ------------------------------------------------------------------
$replacement = "_______________________"
$target = $d.Content
$target.Find.Execute($replacement)

$target.InsertAfter("~")
$bmPos = $target.Start + $replacement.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$marker = $d.Range($bmPos, $bmPos + 1)
$marker.Delete()
